# Auto commit at 2025-09-28  8:36:10.73
#
# Refresh the monthly metrics on the "Metrics" sheet (B2:B13). The
# "today" sheet pulls every one of these numbers in via formulas
# (=Metrics!B2 ... =Metrics!B13, plus the running E/F totals that chain
# off of them), so simply updating the source values here and letting
# Excel recalculate reproduces the rest of the diff automatically -
# including the TODAY()-1 refresh on today!A1.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 402224.54
$metrics.Range("B3").Value  = 324869.39000000007
$metrics.Range("B4").Value  = 127069.59999999999
$metrics.Range("B5").Value  = 15991
$metrics.Range("B6").Value  = 4321475.419999999
$metrics.Range("B7").Value  = 3652396.8699999996
$metrics.Range("B8").Value  = 1256435.2800000003
$metrics.Range("B9").Value  = 167151
$metrics.Range("B10").Value = 32786799.220999829
$metrics.Range("B11").Value = 19682266.940000005
$metrics.Range("B12").Value = 11538144.169999998
$metrics.Range("B13").Value = 1264778

# Restore the recorded selection on Metrics (not the active sheet, so
# select it before moving back to "today" below).
$metrics.Range("E11").Select() | Out-Null

# "today" stays the active/visible tab, with its own recorded selection.
$today = $wb.Worksheets.Item("today")
$today.Range("F8").Select() | Out-Null
